$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the number format (style) of A2 onto A5 so the new time value renders
# the same way as the existing story-point rows.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "This is a different story point on 2"
$ws.Range("D5").Value = "sound/conversation"

$ws.Range("D5").Select() | Out-Null
